$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU per point conversion is handled by PowerPoint COM automatically;
# Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU).

function Set-Pos($shape, [double]$xEmu, [double]$yEmu) {
    $shape.Left = $xEmu / 12700.0
    $shape.Top  = $yEmu / 12700.0
}

# Picture 30 (id 31)
Set-Pos ($s.Shapes.Item("Picture 30"))  1903500 3510169
# TextBox 32 (id 33) - label "left"
Set-Pos ($s.Shapes.Item("TextBox 32"))  1685869 3179278
# Picture 33 (id 34)
Set-Pos ($s.Shapes.Item("Picture 33"))  2074333 4891868
# TextBox 40 (id 41) - label "right"
Set-Pos ($s.Shapes.Item("TextBox 40"))  1821335 4476312
# TextBox 42 (id 43) - "Reference Level"
Set-Pos ($s.Shapes.Item("TextBox 42"))  186266 3191933
# TextBox 43 (id 44) - "Treatment Level"
Set-Pos ($s.Shapes.Item("TextBox 43"))  169333 4817532
# Group 44 (id 45)
Set-Pos ($s.Shapes.Item("Group 44"))  5181598 3143822
# Group 100 (id 101)
Set-Pos ($s.Shapes.Item("Group 100"))  5274731 4828689
# Down Arrow 102 (id 103)
Set-Pos ($s.Shapes.Item("Down Arrow 102"))  6612466 3767667
# Down Arrow 103 (id 104)
Set-Pos ($s.Shapes.Item("Down Arrow 103"))  3539071 4025899
# TextBox 108 (id 109) - beta symbol
Set-Pos ($s.Shapes.Item("TextBox 108"))  6853766 4097867

# Add new TextBox "Biological or Cognitive Functioning" (id 110)
$newBox = $s.Shapes.AddTextbox(1, 8915400 / 12700.0, 457199 / 12700.0, 3115734 / 12700.0, 1200329 / 12700.0)
$newBox.Name = "TextBox 109"
$newBox.TextFrame.WordWrap = $true
$tr = $newBox.TextFrame.TextRange
$tr.Text = "Biological or Cognitive Functioning"
$tr.Font.Size = 24
$tr.Font.Name = "Arial"
